$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 24,13

$data[0,0] = 1.02
$data[0,1] = 1.102879542016455
$data[0,2] = 1.101823131779284
$data[0,3] = 1.115342065623192
$data[0,4] = 1.119365856477834
$data[0,5] = 1
$data[0,6] = $null
$data[0,7] = 1.069860139229108
$data[0,8] = 1.107650339899305
$data[0,9] = 1.10444183671364
$data[0,10] = 1.117927337431595
$data[0,11] = 1.12194134488839
$data[0,12] = 1.109223330431232

$data[1,0] = 1.02
$data[1,1] = 1.104428814165904
$data[1,2] = 1.103113588126888
$data[1,3] = 1.116849896673432
$data[1,4] = 1.120868666987156
$data[1,5] = 1
$data[1,6] = $null
$data[1,7] = 1.070413317494746
$data[1,8] = 1.108866495721032
$data[1,9] = 1.105553382346446
$data[1,10] = 1.119258189275311
$data[1,11] = 1.123267899578378
$data[1,12] = 1.11044121333372

$data[2,0] = 1.02
$data[2,1] = 1.105429908592119
$data[2,2] = 1.10394720798236
$data[2,3] = 1.117824477189137
$data[2,4] = 1.121840017665955
$data[2,5] = 1
$data[2,6] = $null
$data[2,7] = 1.070769129030025
$data[2,8] = 1.109651552711814
$data[2,9] = 1.106270671883338
$data[2,10] = 1.120117726240121
$data[2,11] = 1.124124672558724
$data[2,12] = 1.111227385195477

$data[3,0] = 1.02
$data[3,1] = 1.105850443920678
$data[3,2] = 1.104297334290914
$data[3,3] = 1.118233937361432
$data[3,4] = 1.122248124420556
$data[3,5] = 1
$data[3,6] = $null
$data[3,7] = 1.070918204840487
$data[3,8] = 1.10998114719158
$data[3,9] = 1.106571757788559
$data[3,10] = 1.120478695640382
$data[3,11] = 1.124484483869944
$data[3,12] = 1.111557447737209

$data[4,0] = 1.02
$data[4,1] = 1.105921034827901
$data[4,2] = 1.104356103000197
$data[4,3] = 1.118302672857735
$data[4,4] = 1.122316632927749
$data[4,5] = 1
$data[4,6] = $null
$data[4,7] = 1.070943205660043
$data[4,8] = 1.110036461720492
$data[4,9] = 1.106622284464408
$data[4,10] = 1.12053928189241
$data[4,11] = 1.12454487589967
$data[4,12] = 1.111612840819097

$data[5,0] = 1.02
$data[5,1] = 1.105435529072474
$data[5,2] = 1.10395188766694
$data[5,3] = 1.117829949403611
$data[5,4] = 1.121845471778798
$data[5,5] = 1
$data[5,6] = $null
$data[5,7] = 1.070771122978166
$data[5,8] = 1.109655958504628
$data[5,9] = 1.106274696819055
$data[5,10] = 1.120122551016064
$data[5,11] = 1.124129481845031
$data[5,12] = 1.111231797245022

$data[6,0] = 1.02
$data[6,1] = 1.103403417813521
$data[6,2] = 1.102259539139474
$data[6,3] = 1.11585187333125
$data[6,4] = 1.119873963524597
$data[6,5] = 1
$data[6,6] = $null
$data[6,7] = 1.070047531623855
$data[6,8] = 1.108061737697384
$data[6,9] = 1.10481789676232
$data[6,10] = 1.118377442835421
$data[6,11] = 1.122389994574788
$data[6,12] = 1.109635312461372

$data[7,0] = 1.02
$data[7,1] = 1.099811605339612
$data[7,2] = 1.099266490113693
$data[7,3] = 1.112357622048352
$data[7,4] = 1.116391433604993
$data[7,5] = 1
$data[7,6] = $null
$data[7,7] = 1.068756009494327
$data[7,8] = 1.105237874509316
$data[7,9] = 1.102235616016596
$data[7,10] = 1.115289703204479
$data[7,11] = 1.119312289890864
$data[7,12] = 1.106807439063675

$data[8,0] = 1.02
$data[8,1] = 1.097409208449239
$data[8,2] = 1.097263425646433
$data[8,3] = 1.110021893653887
$data[8,4] = 1.114063624715904
$data[8,5] = 1
$data[8,6] = $null
$data[8,7] = 1.067883737184999
$data[8,8] = 1.103345080912897
$data[8,9] = 1.100503527871204
$data[8,10] = 1.113222319388563
$data[8,11] = 1.11725168821704
$data[8,12] = 1.104911957483243

$data[9,0] = 1.02
$data[9,1] = 1.096366974099771
$data[9,2] = 1.09639416991093
$data[9,3] = 1.109008924947751
$data[9,4] = 1.113054111475331
$data[9,5] = 1
$data[9,6] = $null
$data[9,7] = 1.067503320914576
$data[9,8] = 1.102522973141743
$data[9,9] = 1.099750932005763
$data[9,10] = 1.112324922147226
$data[9,11] = 1.116357250606913
$data[9,12] = 1.104088682224756

$data[10,0] = 1.02
$data[10,1] = 1.095979535009908
$data[10,2] = 1.096070994788715
$data[10,3] = 1.108632417152996
$data[10,4] = 1.112678891248858
$data[10,5] = 1
$data[10,6] = $null
$data[10,7] = 1.067361605498418
$data[10,8] = 1.102217220473305
$data[10,9] = 1.099470988595532
$data[10,10] = 1.111991249662946
$data[10,11] = 1.116024680992449
$data[10,12] = 1.103782495352459

$data[11,0] = 1.02
$data[11,1] = 1.096062656045002
$data[11,2] = 1.096140330405133
$data[11,3] = 1.10871319062421
$data[11,4] = 1.112759388347087
$data[11,5] = 1
$data[11,6] = $null
$data[11,7] = 1.06739202263474
$data[11,8] = 1.102282823009003
$data[11,9] = 1.099531055425028
$data[11,10] = 1.112062839016956
$data[11,11] = 1.116096033615942
$data[11,12] = 1.103848191051283

$data[12,0] = 1.02
$data[12,1] = 1.096334954577222
$data[12,2] = 1.096367462210038
$data[12,3] = 1.108977807764888
$data[12,4] = 1.113023100639926
$data[12,5] = 1
$data[12,6] = $null
$data[12,7] = 1.067491615097904
$data[12,8] = 1.102497707431875
$data[12,9] = 1.099727799908839
$data[12,10] = 1.112297347639767
$data[12,11] = 1.116329767190476
$data[12,12] = 1.104063380634683

$data[13,0] = 1.02
$data[13,1] = 1.096502685751057
$data[13,2] = 1.096507366330872
$data[13,3] = 1.109140814256741
$data[13,4] = 1.113185550165668
$data[13,5] = 1
$data[13,6] = $null
$data[13,7] = 1.067552922611734
$data[13,8] = 1.102630053550178
$data[13,9] = 1.099848968032437
$data[13,10] = 1.112441790964354
$data[13,11] = 1.116473733452766
$data[13,12] = 1.104195914699655

$data[14,0] = 1.02
$data[14,1] = 1.097478335590188
$data[14,2] = 1.097321074327435
$data[14,3] = 1.110089086993645
$data[14,4] = 1.114130589290435
$data[14,5] = 1
$data[14,6] = $null
$data[14,7] = 1.067908926605959
$data[14,8] = 1.10339958787519
$data[14,9] = 1.100553420040706
$data[14,10] = 1.11328182954514
$data[14,11] = 1.117311002442992
$data[14,12] = 1.104966541851675

$data[15,0] = 1.02
$data[15,1] = 1.09808979815348
$data[15,2] = 1.097830974039677
$data[15,3] = 1.110683483798617
$data[15,4] = 1.114722964844645
$data[15,5] = 1
$data[15,6] = $null
$data[15,7] = 1.068131508564245
$data[15,8] = 1.103881617784319
$data[15,9] = 1.100994606014071
$data[15,10] = 1.11380816741321
$data[15,11] = 1.117835609217889
$data[15,12] = 1.105449256298564

$data[16,0] = 1.02
$data[16,1] = 1.098446263448752
$data[16,2] = 1.09812820556598
$data[16,3] = 1.111030033236345
$data[16,4] = 1.115068337848111
$data[16,5] = 1
$data[16,6] = $null
$data[16,7] = 1.068261074995877
$data[16,8] = 1.104162535446532
$data[16,9] = 1.101251692856087
$data[16,10] = 1.114114959024331
$data[16,11] = 1.118141393317058
$data[16,12] = 1.105730572896079

$data[17,0] = 1.02
$data[17,1] = 1.098567776874516
$data[17,2] = 1.098229522778441
$data[17,3] = 1.111148172079602
$data[17,4] = 1.115186075983061
$data[17,5] = 1
$data[17,6] = $null
$data[17,7] = 1.068305209497435
$data[17,8] = 1.104258280243118
$data[17,9] = 1.101339310710567
$data[17,10] = 1.114219531247158
$data[17,11] = 1.118245622375772
$data[17,12] = 1.105826453661259

$data[18,0] = 1.02
$data[18,1] = 1.098024213776016
$data[18,2] = 1.097776285745354
$data[18,3] = 1.110619726445823
$data[18,4] = 1.114659424090756
$data[18,5] = 1
$data[18,6] = $null
$data[18,7] = 1.068107654739012
$data[18,8] = 1.103829925644645
$data[18,9] = 1.10094729678667
$data[18,10] = 1.113751718362691
$data[18,11] = 1.117779345668535
$data[18,12] = 1.105397490750123

$data[19,0] = 1.02
$data[19,1] = 1.096254777991657
$data[19,2] = 1.096300585750393
$data[19,3] = 1.108899891473316
$data[19,4] = 1.112945450690694
$data[19,5] = 1
$data[19,6] = $null
$data[19,7] = 1.06746229899969
$data[19,8] = 1.102434440018446
$data[19,9] = 1.099669874541552
$data[19,10] = 1.11222830012084
$data[19,11] = 1.11626094780577
$data[19,12] = 1.104000023374269

$data[20,0] = 1.02
$data[20,1] = 1.095140483521172
$data[20,2] = 1.095371045223976
$data[20,3] = 1.107817134866009
$data[20,4] = 1.111866402998789
$data[20,5] = 1
$data[20,6] = $null
$data[20,7] = 1.067054153379236
$data[20,8] = 1.101554810302814
$data[20,9] = 1.098864415196457
$data[20,10] = 1.111268500290149
$data[20,11] = 1.115304324959755
$data[20,12] = 1.103119144483547

$data[21,0] = 1.02
$data[21,1] = 1.095731363936584
$data[21,2] = 1.095863976636331
$data[21,3] = 1.108391262783813
$data[21,4] = 1.112438562485433
$data[21,5] = 1
$data[21,6] = $null
$data[21,7] = 1.067270746405319
$data[21,8] = 1.102021332721351
$data[21,9] = 1.099291624039321
$data[21,10] = 1.111777497310426
$data[21,11] = 1.115811635830335
$data[21,12] = 1.103586329417429

$data[22,0] = 1.02
$data[22,1] = 1.098053849121937
$data[22,2] = 1.097800997601421
$data[22,3] = 1.1106485361196
$data[22,4] = 1.114688135885007
$data[22,5] = 1
$data[22,6] = $null
$data[22,7] = 1.068118434065112
$data[22,8] = 1.103853283846197
$data[22,9] = 1.100968674561202
$data[22,10] = 1.113777225915233
$data[22,11] = 1.117804769394275
$data[22,12] = 1.105420882123

$data[23,0] = 1.02
$data[23,1] = 1.100741524293793
$data[23,2] = 1.100041595391321
$data[23,3] = 1.113262035993304
$data[23,4] = 1.117292799165254
$data[23,5] = 1
$data[23,6] = $null
$data[23,7] = 1.069091869165148
$data[23,8] = 1.105969684325254
$data[23,9] = 1.102905035091145
$data[23,10] = 1.116089495851349
$data[23,11] = 1.120109472358442
$data[23,12] = 1.107540288133492

$ws.Range("B2:N25").Value = $data
